# The deck currently ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (used only by the notes master)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet"  (used by the slide
#                            master, i.e. the theme that actually paints
#                            every slide)
#
# The target edit swaps the two theme parts' contents, so the slide master
# (and therefore every slide) switches from the "Integral" palette to the
# plain "Office" palette. Re-apply that palette, color by color, onto the
# presentation's theme color scheme (the PowerPoint object model only
# exposes a single, presentation-wide theme, which is the one bound to the
# slide master / presentation root).

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

# Office theme color scheme, in the fixed ThemeColorScheme index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$cs.Item(1).RGB = 0          # dk1     000000
$cs.Item(2).RGB = 16777215   # lt1     FFFFFF
$cs.Item(3).RGB = 6968388    # dk2     44546A
$cs.Item(4).RGB = 15132391   # lt2     E7E6E6
$cs.Item(5).RGB = 13998939   # accent1 5B9BD5
$cs.Item(6).RGB = 3243501    # accent2 ED7D31
$cs.Item(7).RGB = 10855845   # accent3 A5A5A5
$cs.Item(8).RGB = 49407      # accent4 FFC000
$cs.Item(9).RGB = 12874308   # accent5 4472C4
$cs.Item(10).RGB = 4697456   # accent6 70AD47
$cs.Item(11).RGB = 12673797  # hlink   0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72

# Best-effort: also try to rename the theme / color scheme to match the
# stock "Office Theme" / "Office" naming (harmless no-op on hosts that treat
# these as read-only).
try { $m.Theme.Name = "Office Theme" } catch {}
try { $cs.Name = "Office" } catch {}
